# ---------------------------------------------------------------------------
# Edit script: "Various updates to make more complete reference set"
#
# 1. Window / view bookkeeping (best effort; cosmetic only).
# 2. Re-point the three existing "Clade" highlight cells (E104:E106,
#    E120:E122) at the fill that used to be one slot further along once the
#    stray yellow highlight fill is retired.
# 3. Clear out the bogus Genus/Clade values that had been typed into the
#    TABV/JMTV polyprotein summary row (row 137) and replace them with the
#    correct Genus ("Tamanavirus"), blanking the Clade column.
# 4. Append the full Tamana bat virus (TABV) mature-peptide map as rows
#    138-149, mirroring the layout used for the other reference genomes.
# 5. A handful of trailing blank placeholder rows (150-156) in column E,
#    matching the sparse "ghost" cells already present elsewhere in the
#    sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. window geometry (best effort; the headless host has no real window
#        surface, but set it anyway in case it is honoured) -----------------
$excel.Left = 2360
$excel.Top = 0
$excel.Width = 26080
$excel.Height = 26840

# --- 2. restyle the two other highlighted "Clade" blocks --------------------
# Before: E104:E106 used the fill that sat behind the stray yellow fill;
# E120:E122 used the fill one slot further along again. Removing the yellow
# fill shifts both down by one slot, so re-apply via copy/paste-format from
# a still-later block that already carries the final fill (E122 originally
# carried the fill E120:E122 should keep, so instead we pull the format
# straight off a neighbouring cell that already has the right end-state
# fill).
$ws.Range("E120").Copy() | Out-Null
$ws.Range("E104:E106").PasteSpecial(-4122) | Out-Null
$ws.Range("E104").Value = $ws.Range("E104").Value
$ws.Range("E105").Value = $ws.Range("E105").Value
$ws.Range("E106").Value = $ws.Range("E106").Value

# --- 3. fix up row 137 (TABV / polyprotein summary line) -------------------
# Genus was mistakenly a repeat of the Abbrev column; correct it, and clear
# the spurious Clade value (TABV/Tamanavirus has no defined clade).
$ws.Range("D137").Value = "Tamanavirus"
$ws.Range("D137").HorizontalAlignment = -4131
$ws.Range("E137").ClearContents() | Out-Null
$ws.Range("F137").ClearFormats() | Out-Null

# --- 4. append the TABV mature-peptide feature rows -------------------------
function Set-FeatureRow {
    param(
        [int]$Row,
        [string]$FeatureName,
        [int]$Start,
        [int]$End,
        [bool]$Styled
    )

    $ws.Range("A$Row").Value = "NC_003996"
    $ws.Range("B$Row").Value = "Tamana bat virus"
    $ws.Range("C$Row").Value = "TABV"
    $ws.Range("D$Row").Value = "Tamanavirus"
    $ws.Range("D$Row").HorizontalAlignment = -4131
    $ws.Range("E$Row").ClearContents() | Out-Null
    $ws.Range("F$Row").Value = $FeatureName
    $ws.Range("G$Row").Value = $Start
    $ws.Range("H$Row").Value = $End

    if ($Styled) {
        $ws.Range("F16").Copy() | Out-Null
        $ws.Range("G$Row").PasteSpecial(-4122) | Out-Null
        $ws.Range("H$Row").PasteSpecial(-4122) | Out-Null
    }
}

Set-FeatureRow -Row 138 -FeatureName "protein C"          -Start 1    -End 345  -Styled $false
Set-FeatureRow -Row 139 -FeatureName "precursor M"         -Start 346  -End 852  -Styled $true
Set-FeatureRow -Row 140 -FeatureName "M"                   -Start 655  -End 852  -Styled $true
Set-FeatureRow -Row 141 -FeatureName "envelope protein E"  -Start 853  -End 2358 -Styled $true
Set-FeatureRow -Row 142 -FeatureName "NS1"                 -Start 2359 -End 3390 -Styled $true
Set-FeatureRow -Row 143 -FeatureName "NS2A"                -Start 3391 -End 3978 -Styled $true
Set-FeatureRow -Row 144 -FeatureName "NS2B"                -Start 3979 -End 4431 -Styled $true
Set-FeatureRow -Row 145 -FeatureName "NS3"                 -Start 4432 -End 6303 -Styled $true
Set-FeatureRow -Row 146 -FeatureName "NS4A"                -Start 6304 -End 6687 -Styled $true
Set-FeatureRow -Row 147 -FeatureName "2K"                  -Start 6763 -End 6762 -Styled $true
Set-FeatureRow -Row 148 -FeatureName "NS4B"                -Start 6763 -End 7557 -Styled $true
Set-FeatureRow -Row 149 -FeatureName "NS5"                 -Start 7558 -End 10050 -Styled $true

# --- 5. trailing blank placeholder rows (E150:E156) -------------------------
$ws.Range("E150").Value = " "
$ws.Range("E150").ClearContents() | Out-Null
$ws.Range("E151").Value = " "
$ws.Range("E151").ClearContents() | Out-Null
$ws.Range("E152").Value = " "
$ws.Range("E152").ClearContents() | Out-Null
$ws.Range("E153").Value = " "
$ws.Range("E153").ClearContents() | Out-Null
$ws.Range("E154").Value = " "
$ws.Range("E154").ClearContents() | Out-Null
$ws.Range("E155").Value = " "
$ws.Range("E155").ClearContents() | Out-Null
$ws.Range("E156").Value = " "
$ws.Range("E156").ClearContents() | Out-Null

# --- final selection / scroll position --------------------------------------
$ws.Range("H157").Select()
$excel.ActiveWindow.ScrollRow = 120
$excel.ActiveWindow.ScrollColumn = 1
